$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.360.42'
$ws.Range('E2').Value = '  +2.33%  '
$ws.Range('D3').Value = '1.662.20'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.20'
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('E8').Value = '  +1.41%  '
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.01'
$ws.Range('E10').Value = '  +4.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0849'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').Value = '1.895.25'
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('D13').Value = '1.656.89'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.22'
$ws.Range('E14').Value = '  +1.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.535'
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.36'
$ws.Range('E16').Value = '  +4.18%  '
$ws.Range('D17').Value = '27.340.42'
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '222.68'
$ws.Range('E19').Value = '  +3.62%  '
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('E21').Value = '  +9.27%  '
$ws.Range('E22').Value = '  +2.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.51'
$ws.Range('E23').Value = '  +5.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.30'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.28'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  +4.12%  '
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.07'
$ws.Range('E29').Value = '  +3.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0516'
$ws.Range('E30').Value = '  +1.63%  '
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.40'
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('E34').Value = '  +2.36%  '
$ws.Range('D35').Value = '1.265.79'
$ws.Range('E35').Value = '  -1.64%  '
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.538'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.836'
$ws.Range('E39').Value = '  +2.44%  '
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.816'
$ws.Range('E41').Value = '  +1.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.41'
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('D43').Value = '1.806.96'
$ws.Range('E43').Value = '  +1.49%  '
$ws.Range('E44').Value = '  -3.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.87'
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.61'
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('E47').Value = '  +1.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0520'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('E49').Value = '  +2.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.67'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('E51').Value = '  +0.37%  '
